$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("策略更新")

# Change the 日/夜盘 value from 夜盘 to 日盘
$ws.Range("C3").Value = "日盘"

# Change the description from the old maintenance note to the new 锁仓测试 text
$ws.Range("E3").Value = "锁仓测试"

# Move the date back one day (2017-02-22 -> 2017-02-21)
$ws.Range("D3").Value = (Get-Date -Year 2017 -Month 2 -Day 21 -Hour 0 -Minute 0 -Second 0)

# Update the active selection to F4
$ws.Range("F4").Select()
